$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1777.625
$ws.Range("J39").Value = 1262.3334
$ws.Range("L39").Value = 3787.0002
$ws.Range("N39").Value = -4379.0002
$ws.Range("H92").Value = 637.7
$ws.Range("J92").Value = 662
$ws.Range("L92").Value = 662
$ws.Range("N92").Value = -3158
$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 11225.714
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 11225.714
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 33677.142
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -38737.142
$ws.Range("H137").Value = 3321.1875
$ws.Range("I137").Value = 1558.2727
$ws.Range("J137").Value = 7199.6
$ws.Range("K137").Value = 4674.8181
$ws.Range("L137").Value = 21598.8
$ws.Range("M137").Value = -2124.8181
$ws.Range("N137").Value = -26698.8
$ws.Range("H138").Value = 2144.3618
$ws.Range("I138").Value = 1515.4
$ws.Range("J138").Value = 2610.2593
$ws.Range("K138").Value = 4546.200000000001
$ws.Range("L138").Value = 7830.777900000001
$ws.Range("M138").Value = 593.7999999999993
$ws.Range("N138").Value = -18110.7779
$ws.Range("H141").Value = 1254.8
$ws.Range("I141").Value = 1254.8
$ws.Range("K141").Value = 3764.4
$ws.Range("M141").Value = 1415.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4333.541
$ws.Range("I32").Value = 3935.9824
$ws.Range("K32").Value = 3935.9824
$ws.Range("M32").Value = -3648.9824
$ws.Range("H45").Value = 8407
$ws.Range("I45").Value = 14688.375
$ws.Range("K45").Value = 14688.375
$ws.Range("M45").Value = -14311.375
$ws.Range("H110").Value = 3957.1667
$ws.Range("I110").Value = 1955.4286
$ws.Range("K110").Value = 1955.4286
$ws.Range("M110").Value = 89.57140000000004
$ws.Range("H122").Value = 3703.6316
$ws.Range("I122").Value = 3332.111
$ws.Range("K122").Value = 9996.332999999999
$ws.Range("M122").Value = -7546.332999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 49593
$ws.Range("I107").Value = 60770.824
$ws.Range("J107").Value = 2087.25
$ws.Range("K107").Value = 60770.824
$ws.Range("L107").Value = 2087.25
$ws.Range("M107").Value = -58850.824
$ws.Range("N107").Value = -5927.25
$ws.Range("H134").Value = 1286.9714
$ws.Range("I134").Value = 1229.0625
$ws.Range("K134").Value = 3687.1875
$ws.Range("M134").Value = -1152.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2339.5
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 2452.6667
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2452.6667
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = -3026.6667
$ws.Range("H25").Value = 5000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 5000
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -5348
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H60").Value = 14237.5
$ws.Range("I60").Value = 8900
$ws.Range("K60").Value = 8900
$ws.Range("M60").Value = -8389
$ws.Range("H68").Value = 20995
$ws.Range("J68").Value = 20995
$ws.Range("L68").Value = 20995
$ws.Range("N68").Value = -22493
$ws.Range("H71").Value = 20995
$ws.Range("J71").Value = 20995
$ws.Range("L71").Value = 62985
$ws.Range("N71").Value = -70473
$ws.Range("H99").Value = 3577.818
$ws.Range("I99").Value = 2411.5833
$ws.Range("K99").Value = 2411.5833
$ws.Range("M99").Value = -913.5832999999998
$ws.Range("H106").Value = 40671
$ws.Range("J106").Value = 40671
$ws.Range("L106").Value = 40671
$ws.Range("N106").Value = -43195
$ws.Range("H107").Value = 2338.4
$ws.Range("I107").Value = 1961.4445
$ws.Range("K107").Value = 1961.4445
$ws.Range("M107").Value = -41.44450000000006
$ws.Range("H113").Value = 2339.5
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 2452.6667
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 2452.6667
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -6792.6667
$ws.Range("H122").Value = 47442.26
$ws.Range("I122").Value = 78710.53999999999
$ws.Range("K122").Value = 236131.62
$ws.Range("M122").Value = -233681.62
$ws.Range("H126").Value = 3577.818
$ws.Range("I126").Value = 2411.5833
$ws.Range("K126").Value = 7234.749899999999
$ws.Range("M126").Value = -4764.749899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3147.7
$ws.Range("I5").Value = 1019.6
$ws.Range("K5").Value = 3058.8
$ws.Range("M5").Value = -2946.8
$ws.Range("H39").Value = 75090.86
$ws.Range("J39").Value = 4633.4546
$ws.Range("L39").Value = 13900.3638
$ws.Range("N39").Value = -14488.3638
$ws.Range("H131").Value = 54612.684
$ws.Range("I131").Value = 144624.28
$ws.Range("J131").Value = 2105.9167
$ws.Range("K131").Value = 433872.84
$ws.Range("L131").Value = 6317.750100000001
$ws.Range("M131").Value = -428832.84
$ws.Range("N131").Value = -16397.7501
$ws.Range("H135").Value = 3147.7
$ws.Range("I135").Value = 1019.6
$ws.Range("K135").Value = 9176.4
$ws.Range("M135").Value = -6641.4
$ws.Range("H137").Value = 4119.846
$ws.Range("J137").Value = 4348.6665
$ws.Range("L137").Value = 13045.9995
$ws.Range("N137").Value = -23245.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2910.6
$ws.Range("I113").Value = 2449.7144
$ws.Range("J113").Value = 3313.875
$ws.Range("K113").Value = 2449.7144
$ws.Range("L113").Value = 3313.875
$ws.Range("M113").Value = -279.7143999999998
$ws.Range("N113").Value = -7653.875
$ws.Range("H122").Value = 3098.76
$ws.Range("I122").Value = 2972.9333
$ws.Range("K122").Value = 8918.7999
$ws.Range("M122").Value = -6468.7999
$ws.Range("H126").Value = 3084
$ws.Range("I126").Value = 3084
$ws.Range("K126").Value = 9252
$ws.Range("M126").Value = -6782
$ws.Range("H132").Value = 4995.5938
$ws.Range("I132").Value = 3611.7144
$ws.Range("K132").Value = 10835.1432
$ws.Range("M132").Value = -8305.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1666
$ws.Range("I22").Value = 1499
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1499
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1204
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1666
$ws.Range("I27").Value = 1499
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1499
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1392
$ws.Range("N27").Value = -2214
$ws.Range("H46").Value = 37176.082
$ws.Range("I46").Value = 72119.664
$ws.Range("K46").Value = 72119.664
$ws.Range("M46").Value = -71931.664
$ws.Range("H61").Value = 106061.375
$ws.Range("I61").Value = 92101.63
$ws.Range("K61").Value = 92101.63
$ws.Range("M61").Value = -91899.63
$ws.Range("H113").Value = 106061.375
$ws.Range("I113").Value = 92101.63
$ws.Range("K113").Value = 92101.63
$ws.Range("M113").Value = -89931.63
$ws.Range("H132").Value = 3842.6047
$ws.Range("I132").Value = 3233.9119
$ws.Range("J132").Value = 6142.1113
$ws.Range("K132").Value = 9701.735700000001
$ws.Range("L132").Value = 18426.3339
$ws.Range("M132").Value = -7171.735700000001
$ws.Range("N132").Value = -23486.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1293.7455
$ws.Range("I122").Value = 1178.7174
$ws.Range("K122").Value = 3536.1522
$ws.Range("M122").Value = -1086.1522
$ws.Range("H136").Value = 1193.4706
$ws.Range("I136").Value = 922.3333
$ws.Range("K136").Value = 2766.9999
$ws.Range("M136").Value = -216.9998999999998
